$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the very top; this pushes all existing
# data (old rows 1-31) down to rows 3-33.
$ws.Rows("1:2").Insert()

# Populate the new row 1 with a simple numeric index header (0..10).
for ($i = 0; $i -lt 11; $i++) {
    $col = [char](65 + $i)
    $ws.Range("$col`1").Value = $i
}

# New row 2 only contains a single label in column E.
$ws.Range("E2").Value = "Drive"

# Give the new row 1 the same (bold / centered / bordered) formatting
# that the original header row used, by copying the formats from the
# shifted-down header (now row 3, which still carries that formatting).
$ws.Range("A3:K3").Copy()
$ws.Range("A1:K1").PasteSpecial(-4122)

# The old header row (now row 3) loses its special formatting.
$ws.Range("A3:K3").ClearFormats()

# The old header's J/K labels ("thread_size" / "material_surface") are
# removed, leaving those two cells blank.
$ws.Range("J3").ClearContents()
$ws.Range("K3").ClearContents()
